# Updated with data from Apr 10
# - adds the missing sd_total (column B) value for 4/9 (row 36)
# - appends a brand-new row (row 37) of data for 4/10/2020 (date serial 43931)
# - moves the selection/view to the new last row/column, matching Excel's
#   behaviour after typing data into a new row at the bottom-right of the
#   used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-missing sd_total figure for 4/9 (row 36) ---
$ws.Range("B36").Value = 1255

# --- Append the new row for 4/10/2020 (Excel serial date 43931) ---
$ws.Range("A37").Value = 43931

$ws.Range("C37").Value = 1693
$ws.Range("D37").Value = 10
$ws.Range("F37").Value = 21
$ws.Range("H37").Value = 264
$ws.Range("J37").Value = 335
$ws.Range("L37").Value = 292
$ws.Range("N37").Value = 305
$ws.Range("P37").Value = 235
$ws.Range("R37").Value = 137
$ws.Range("T37").Value = 92
$ws.Range("V37").Value = 2
$ws.Range("X37").Value = 821
$ws.Range("Y37").Value = 867
$ws.Range("Z37").Value = 5
$ws.Range("AA37").Value = 374
$ws.Range("AB37").Value = 140
$ws.Range("AC37").Value = 44
$ws.Range("AD37").Value = 48
$ws.Range("AE37").Value = 163
$ws.Range("AF37").Value = 5
$ws.Range("AG37").Value = 8
$ws.Range("AH37").Value = 97
$ws.Range("AI37").Value = 33
$ws.Range("AJ37").Value = 46
$ws.Range("AK37").Value = 10
$ws.Range("AL37").Value = 33
$ws.Range("AM37").Value = 14
$ws.Range("AN37").Value = 33
$ws.Range("AO37").Value = 36
$ws.Range("AP37").Value = 17
$ws.Range("AQ37").Value = 849
$ws.Range("AR37").Value = 17
$ws.Range("AS37").Value = 17
$ws.Range("AT37").Value = 5
$ws.Range("AU37").Value = 28
$ws.Range("AV37").Value = 1
$ws.Range("AW37").Value = 13
$ws.Range("AX37").Value = 1
$ws.Range("AY37").Value = 1
$ws.Range("AZ37").Value = 8
$ws.Range("BA37").Value = 4
$ws.Range("BB37").Value = 14
$ws.Range("BC37").Value = 2
$ws.Range("BD37").Value = 10
$ws.Range("BE37").Value = 14
$ws.Range("BG37").Value = 42
$ws.Range("BH37").Value = 3
$ws.Range("BI37").Value = 102

# --- Match the saved selection / scroll position from the source workbook ---
$excel.ActiveWindow.ScrollColumn = 52
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BI37").Select()
